$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '72.170.63'
$ws.Range("E2").Value = '  +0.28%  '

# Row 3
$ws.Range("D3").Value = '4.028.63'
$ws.Range("E3").Value = '  -0.47%  '

# Row 4
$ws.Range("E4").Value = '  +0.14%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '532.35'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.47%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '151.18'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.54%  '

# Row 7
$ws.Range("E7").Value = '  +12.25%  '

# Row 8
$ws.Range("E8").Value = '  +0.07%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.753'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.58%  '

# Row 10
$ws.Range("E10").Value = '  -2.97%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000327'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.34%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '47.82'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.96%  '

# Row 13
$ws.Range("D13").Value = '4.689.10'
$ws.Range("E13").Value = '  +0.35%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.67'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.23%  '

# Row 15
$ws.Range("D15").Value = '4.035.42'
$ws.Range("E15").Value = '  +0.11%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.09'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.97%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '20.58'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.46%  '

# Row 18
$ws.Range("E18").Value = '  -0.57%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.19'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.74%  '

# Row 20
$ws.Range("D20").Value = '72.170.07'
$ws.Range("E20").Value = '  +0.47%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '431.36'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.66%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '98.35'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.82%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.47'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.10%  '

# Row 24
$ws.Range("E24").Value = '  +4.51%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '14.17'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.50%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.10'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -10.18%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.72'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.84%  '

# Row 28
$ws.Range("E28").Value = '  +1.24%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '36.78'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.36%  '

# Row 30
$ws.Range("E30").Value = '  +22.57%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '13.36'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.97%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.128'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.21%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '676.14'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.55%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.04'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.67%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '44.60'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +8.96%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '65.96'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.70%  '

# Row 37
$ws.Range("E37").Value = '  -0.13%  '

# Row 38
$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.152'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.26%  '

# Row 39
$ws.Range("B39").Value = 'PEPE'
$ws.Range("C39").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D39").Value = '0.0₃0830'
$ws.Range("E39").Value = '  -9.22%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.38'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.70%  '

# Row 41
$ws.Range("E41").Value = '  -0.28%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.10%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0485'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.93%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.19'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.04%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.150'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.18%  '

# Row 46
$ws.Range("B46").Value = 'THORChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.75'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.77%  '

# Row 47
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.41'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.72%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.62'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -6.77%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.01'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -6.65%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.000271'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.74%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '144.84'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.58%  '
